$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the second header row ("Hiver"/"Eté"/"Année" labels). This shifts
#    all data rows (old rows 3-15) up by one (new rows 2-14).
$ws.Rows("2:2").Delete()

# 2. Clear the remaining stray header strings/styles left in row 1
#    (E1,G1,I1,J1,K1 carried the old "mation"/"pompes)"/"Hiver"/"Eté"/"Année"
#    leftovers) so we can rebuild the header cleanly from scratch.
$ws.Range("A1:K1").ClearFormats()
$ws.Range("A1:K1").ClearContents()

# 3. Write the new header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# 4. Apply a plain (non number-formatted) copy of the existing data font to the
#    numeric header cells F1:K1, matching the new style added to the
#    workbook's cellXfs (font 1, no explicit number format).
$headerStyle = $wb.Styles.Add("NE2014HeaderStyle")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "NE2014HeaderStyle"
$wb.Styles.Item("NE2014HeaderStyle").Delete()

# 5. Restore a tidy view: select the first data row instead of the old
#    scrolled/selected state.
$ws.Activate()
$ws.Range("A2:K2").Select()
